$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-36 down to 31-37
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price record
$ws.Cells.Item(30,1).Value  = 4
$ws.Cells.Item(30,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30,3).Value  = "Los Lagos"
$ws.Cells.Item(30,4).Value  = 44449
$ws.Cells.Item(30,5).Value  = 10
$ws.Cells.Item(30,6).Value  = 100112026
$ws.Cells.Item(30,7).Value  = "Haba"
$ws.Cells.Item(30,8).Value  = "Sin especificar"
$ws.Cells.Item(30,9).Value  = "Primera"
$ws.Cells.Item(30,10).Value = 80
$ws.Cells.Item(30,11).Value = 17000
$ws.Cells.Item(30,12).Value = 17000
$ws.Cells.Item(30,13).Value = 17000
$ws.Cells.Item(30,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(30,15).Value = "Provincia de Limarí"
$ws.Cells.Item(30,16).Value = 680
$ws.Cells.Item(30,17).Value = 25
$ws.Cells.Item(30,18).Value = "Hortaliza"
